# "update all missing trips" - append 3 new daily rows (4/20/16 - 4/22/16)
# to the bottom of the Data table (rows 37-39), matching the pattern of the
# existing trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# -- Row 37 (2016-04-20 / serial 42510) -------------------------------
$ws.Range("A37").Value = 42510
$ws.Range("B37").Value = 138
$ws.Range("C37").Value = 130
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 7
$ws.Range("F37").Value = 131
$ws.Range("G37").Value = 0.94927536231884058
$ws.Range("H37").Value = 44.964734298409894
$ws.Range("I37").Value = 34.516666667768732
$ws.Range("J37").Value = 63.233333331299946

# -- Row 38 (2016-04-21 / serial 42511) -------------------------------
$ws.Range("A38").Value = 42511
$ws.Range("B38").Value = 139
$ws.Range("C38").Value = 125
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 14
$ws.Range("F38").Value = 125
$ws.Range("G38").Value = 0.89928057553956831
$ws.Range("H38").Value = 45.517146282958322
$ws.Range("I38").Value = 36.033333332743496
$ws.Range("J38").Value = 126.35000000009313

# -- Row 39 (2016-04-22 / serial 42512) -------------------------------
$ws.Range("A39").Value = 42512
$ws.Range("B39").Value = 130
$ws.Range("C39").Value = 121
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 9
$ws.Range("F39").Value = 121
$ws.Range("G39").Value = 0.93076923076923079
$ws.Range("H39").Value = 44.905000000228533
$ws.Range("I39").Value = 35.216666663764045
$ws.Range("J39").Value = 65.900000003166497

# Column A keeps the bordered date format used by the rows directly above
# it (s="12") rather than the worksheet's default/column date style, so
# copy that formatting down instead of leaving the freshly-typed default.
$ws.Range("A36").Copy() | Out-Null
$ws.Range("A37:A39").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Leave the selection where the author left it after entering the last
# row of new data.
$ws.Range("H39:J39").Select() | Out-Null
